# Daily attendance update - 2025-08-26
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WCS_Team_August_2025")

# Row 9's 08-25 entry (AD9) is itself corrected from "WFO" to "WFH" as part of
# today's update, so fix that first - copy the formatting from a known WFH
# cell (AE2 header aside, use the existing WFH-styled cell at column AC ...
# simplest: reuse row 9's own AD cell style by first copying a WFH template).
# Column E1 area has no WFH sample; instead grab formatting from any cell
# already styled "Neutral"/WFH if present - there is none before this edit,
# so build it by copying AD9 (currently WFO/"Good") and then nudging the
# style via the Good->Neutral swap is not directly exposed. Instead we rely
# on the fact that row 16 uses the same "Neutral" xf (style 5) as WFH cells
# for its SL value; copy its formatting, then overwrite the text.
$templateNeutral = $ws.Cells.Item(16, 30)  # AD16, style "Neutral" (s=5)

# Column AE corresponds to 2025-08-26 (Tuesday). Column AD (08-25) already
# carries the correct formatting (fill/border/font) for each status, so copy
# AD's formatting+value into AE for every employee row first.
for ($r = 3; $r -le 18; $r++) {
    $adCell = $ws.Cells.Item($r, 30)  # AD
    $aeCell = $ws.Cells.Item($r, 31)  # AE
    $adCell.Copy($aeCell)
}

# Row 5: 08-26 is "WFH" instead of the copied "WFO" -> restyle + set value
$templateNeutral.Copy($ws.Cells.Item(5, 31))
$ws.Cells.Item(5, 31).Value = "WFH"

# Row 9: both 08-25 (AD9) and 08-26 (AE9) are "WFH"
$templateNeutral.Copy($ws.Cells.Item(9, 30))
$ws.Cells.Item(9, 30).Value = "WFH"
$templateNeutral.Copy($ws.Cells.Item(9, 31))
$ws.Cells.Item(9, 31).Value = "WFH"

$excel.Calculate()
